$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append three new data rows (75-77) to the temperature log table
$ws.Range("A75").Value = 20250702
$ws.Range("B75").Value = "plate43"
$ws.Range("C75").Value = "T1"
$ws.Range("D75").Value = 32.5

$ws.Range("A76").Value = 20250702
$ws.Range("B76").Value = "plate45"
$ws.Range("C76").Value = "T1"
$ws.Range("D76").Value = 34.5

$ws.Range("A77").Value = 20250702
$ws.Range("B77").Value = "plate47"
$ws.Range("C77").Value = "T1"
$ws.Range("D77").Value = 33.5

# Update the selected cell to match the saved view state
$ws.Range("B77").Select()
